$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C24 currently uses numFmtId 167 (YYYY-MM-DD); change it to numFmtId 165 (YYYY-MM-DD HH:MM:SS)
$ws.Range("C24").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add new row 25 values
$ws.Cells.Item(25, 1).Value = 776.14
$ws.Cells.Item(25, 2).Value = 673.0359999999999
$ws.Cells.Item(25, 3).Value = 45756
$ws.Range("C25").NumberFormat = "YYYY-MM-DD"
